$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties in AD1:AF1, matching the bold/bordered
# header formatting already used by the rest of row 1 (e.g. AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$header = $ws.Range("AD1:AF1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data rows 2-46: same team record (Wins/Losses/Ties) applied to every player row
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 90  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 72  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
